$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (leading apostrophe forces text storage, matching the
# original inline-string/text semantics for values that look numeric/date-like)
$ws.Range("A2").Value = "'381943"
$ws.Range("B2").Value = "iliya"
$ws.Range("D2").Value = "'0441201423"
$ws.Range("F2").Value = "'02/27/2023"

# Row 3 is removed entirely, shrinking the used range to A1:F2
$ws.Rows.Item(3).Delete()
